# Updates crypto price/volume data per Aug 21 2023 GitHub Actions refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.225.65"
$ws.Range("E2").Value = "  -0.20%  "
$ws.Range("D3").Value = "1.682.81"
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.59"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5272"
$ws.Range("E6").Value = "  -0.26%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2699"
$ws.Range("E8").Value = "  +0.41%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06382"
$ws.Range("E9").Value = "  -1.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.57"
$ws.Range("E10").Value = "  -1.58%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07658"
$ws.Range("E11").Value = "  +1.94%  "
$ws.Range("D12").Value = "1.695.32"
$ws.Range("E12").Value = "  -0.19%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.510"
$ws.Range("E13").Value = "  -0.08%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5761"
$ws.Range("E14").Value = "  -0.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008335"
$ws.Range("E15").Value = "  -1.87%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.16"
$ws.Range("E16").Value = "  +2.31%  "
$ws.Range("D17").Value = "26.260.21"
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.873"
$ws.Range("E19").Value = "  -0.88%  "
$ws.Range("E20").Value = "  -0.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "189.59"
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.230"
$ws.Range("E22").Value = "  +0.56%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.009"
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "148.83"
$ws.Range("E24").Value = "  +2.83%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.801"
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1260"
$ws.Range("E26").Value = "  -1.29%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.74"
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06268"
$ws.Range("E28").Value = "  -3.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.371"
$ws.Range("E29").Value = "  +0.62%  "
$ws.Range("E30").Value = "  +0.16%  "
$ws.Range("E31").Value = "  -0.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.563"
$ws.Range("E32").Value = "  -0.63%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.687"
$ws.Range("E33").Value = "  +2.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.024"
$ws.Range("E34").Value = "  -0.48%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6128"
$ws.Range("E35").Value = "  -0.94%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.421"
$ws.Range("E36").Value = "  +0.68%  "
$ws.Range("E37").Value = "  +0.94%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.182"
$ws.Range("E38").Value = "  -1.56%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.8980"
$ws.Range("E39").Value = "  +2.98%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01622"
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "1.100.38"
$ws.Range("E41").Value = "  -1.42%  "
$ws.Range("E42").Value = "  -0.34%  "
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("D44").Value = "1.834.13"
$ws.Range("E44").Value = "  +0.45%  "
$ws.Range("E45").Value = "  -1.36%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.42"
$ws.Range("E46").Value = "  +0.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.005"
$ws.Range("E47").Value = "  -0.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.067"
$ws.Range("E48").Value = "  -1.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05276"
$ws.Range("E49").Value = "  +0.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4289"
$ws.Range("E50").Value = "  -0.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.021"
$ws.Range("E51").Value = "  -0.42%  "
